$wb  = $excel.ActiveWorkbook
$ws0 = $wb.Worksheets.Item(1)

# --- Add a second worksheet ("Feuil1") right after "Sheet_0" -----------------
$newSheet = $wb.Worksheets.Add($null, $ws0, $null, $null)
$newSheet.Name = "Feuil1"

# --- Row 1: duplicate the header formulas together with their bold/centered
#     style (copy/paste-special keeps the existing style index instead of
#     creating a brand new one). -4122 == xlPasteFormats.
$ws0.Range("A1:O1").Copy()
$newSheet.Range("A1:O1").PasteSpecial(-4122)
for ($col = 1; $col -le 15; $col++) {
    $srcCell = $ws0.Cells.Item(1, $col)
    $dstCell = $newSheet.Cells.Item(1, $col)
    $dstCell.Formula = $srcCell.Formula
}

# --- Row 2: duplicate the literal/shared-string values -----------------------
for ($col = 1; $col -le 15; $col++) {
    $srcCell = $ws0.Cells.Item(2, $col)
    $dstCell = $newSheet.Cells.Item(2, $col)
    $dstCell.Value2 = $srcCell.Value2
}

# --- Selections: "Sheet_0" keeps a whole-used-rows selection and is no longer
#     the active tab; "Feuil1" becomes active with E9 selected. -------------
$newSheet.Range("E9").Select()
$ws0.Range("A1:XFD2").Select()
$newSheet.Select()
